$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header renames (Инд# -> "И #", Лаб# -> "Л #", Лаб3-4 -> "Л 3 4") ---
$ws.Range("D1").Value  = "И 1"
$ws.Range("E1").Value  = "И 2"
$ws.Range("F1").Value  = "И 3"
$ws.Range("G1").Value  = "И 4"
$ws.Range("H1").Value  = "И 5"
$ws.Range("I1").Value  = "И 6"
$ws.Range("J1").Value  = "И 7"
$ws.Range("K1").Value  = "И 8"
$ws.Range("L1").Value  = "И 9"
$ws.Range("M1").Value  = "И 10"
$ws.Range("N1").Value  = "И 11"
$ws.Range("O1").Value  = "И 12"
$ws.Range("Q1").Value  = "Л 1"
$ws.Range("R1").Value  = "Л 2"
$ws.Range("S1").Value  = "Л 3 4"
$ws.Range("T1").Value  = "Л 5"
$ws.Range("U1").Value  = "Л 6"
$ws.Range("V1").Value  = "Л 7"
$ws.Range("W1").Value  = "Л 8"
$ws.Range("X1").Value  = "Л 9"
$ws.Range("Y1").Value  = "Л 10"
$ws.Range("Z1").Value  = "Л 11"
# P1 ("Тест"), AD1 ("Вариант"), AE1 ("оценко") are unchanged.

# --- Column widths got much narrower ---
# (ColumnWidth is quantized to 1/6 of a character by the host, so we pick the
# value whose round-trip lands closest to the authored OOXML width.)
$ws.Range("D1:L1").EntireColumn.ColumnWidth  = 1.5
$ws.Range("M1:O1").EntireColumn.ColumnWidth  = 2.5
$ws.Range("P1").EntireColumn.ColumnWidth     = 4
$ws.Range("Q1:R1").EntireColumn.ColumnWidth  = 1.6666666666666667
$ws.Range("S1").EntireColumn.ColumnWidth     = 2.5
$ws.Range("T1:W1").EntireColumn.ColumnWidth  = 1.1666666666666667
$ws.Range("X1").EntireColumn.ColumnWidth     = 1.3333333333333333
$ws.Range("Y1:Z1").EntireColumn.ColumnWidth  = 2.3333333333333335

# --- View: zoom 115% -> 130%, scroll/selection moved ---
$excel.ActiveWindow.Zoom = 130
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("AA1").Select()
